$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.951.32'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.846.17'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.64%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.38'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4683'
$ws.Range('E7').Value = '  +3.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3660'
$ws.Range('E8').Value = '  +1.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07151'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9274'
$ws.Range('E10').Value = '  +3.45%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.58'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07699'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.914.87'
$ws.Range('E13').Value = '  +5.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.287'
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.399'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.22'
$ws.Range('E16').Value = '  +3.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008620'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.981.99'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.42'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.027'
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.940'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.40'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.23'
$ws.Range('E26').Value = '  +2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.022'
$ws.Range('E27').Value = '  -2.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.27'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.885'
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.222'
$ws.Range('E31').Value = '  +3.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.179'
$ws.Range('E32').Value = '  +5.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7467'
$ws.Range('E33').Value = '  -0.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.795'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.469'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.084'
$ws.Range('E36').Value = '  +1.19%  '
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.972'
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05191'
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5206'
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.910'
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1519'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.137'
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.54'
$ws.Range('E44').Value = '  +5.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4693'
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.46'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '65.12'
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06034'
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8890'
$ws.Range('E51').Value = '  +5.09%  '
